# Scheduled-runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve
# sheets. Values below are hard-coded snapshots (no formulas backing
# these cells), matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value = 1750
$ws.Range("L17").Value = 5250
$ws.Range("N17").Value = -5586
$ws.Range("H113").Value = 3731.3333
$ws.Range("I113").Value = 4399.3335
$ws.Range("J113").Value = 3063.3333
$ws.Range("K113").Value = 4399.3335
$ws.Range("L113").Value = 3063.3333
$ws.Range("M113").Value = -1145.3335
$ws.Range("N113").Value = -9571.3333
$ws.Range("H137").Value = 2355.1177
$ws.Range("J137").Value = 3166.5
$ws.Range("L137").Value = 9499.5
$ws.Range("N137").Value = -14599.5
$ws.Range("H138").Value = 2308.9666
$ws.Range("I138").Value = 1131.8572
$ws.Range("K138").Value = 3395.5716
$ws.Range("M138").Value = 1744.4284
$ws.Range("H140").Value = 90081.836
$ws.Range("J140").Value = 90081.836
$ws.Range("L140").Value = 90081.836
$ws.Range("N140").Value = -100441.836

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 76254.914
$ws.Range("I74").Value = 42018.363
$ws.Range("J74").Value = 452857
$ws.Range("K74").Value = 42018.363
$ws.Range("L74").Value = 452857
$ws.Range("M74").Value = -41144.363
$ws.Range("N74").Value = -454605
$ws.Range("H77").Value = 76254.914
$ws.Range("I77").Value = 42018.363
$ws.Range("J77").Value = 452857
$ws.Range("K77").Value = 210091.815
$ws.Range("L77").Value = 2264285
$ws.Range("M77").Value = -205723.815
$ws.Range("N77").Value = -2273021
$ws.Range("H132").Value = 7510.393
$ws.Range("I132").Value = 7553.519
$ws.Range("J132").Value = 6949.75
$ws.Range("K132").Value = 22660.557
$ws.Range("L132").Value = 20849.25
$ws.Range("M132").Value = -20130.557
$ws.Range("N132").Value = -25909.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 49542.668
$ws.Range("J76").Value = 49542.668
$ws.Range("L76").Value = 49542.668
$ws.Range("N76").Value = -50172.668
$ws.Range("H79").Value = 49542.668
$ws.Range("J79").Value = 49542.668
$ws.Range("L79").Value = 49542.668
$ws.Range("N79").Value = -51726.668
$ws.Range("H134").Value = 2564.3928
$ws.Range("I134").Value = 1783.3829
$ws.Range("J134").Value = 6643
$ws.Range("K134").Value = 5350.1487
$ws.Range("L134").Value = 19929
$ws.Range("M134").Value = -2815.1487
$ws.Range("N134").Value = -24999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 333.625
$ws.Range("I7").Value = 172.25
$ws.Range("J7").Value = 495
$ws.Range("K7").Value = 172.25
$ws.Range("L7").Value = 495
$ws.Range("M7").Value = -59.25
$ws.Range("N7").Value = -721
$ws.Range("H31").Value = 2735.3333
$ws.Range("I31").Value = 2281
$ws.Range("K31").Value = 2281
$ws.Range("M31").Value = -1986
$ws.Range("H34").Value = 2735.3333
$ws.Range("I34").Value = 2281
$ws.Range("K34").Value = 2281
$ws.Range("M34").Value = -2079
$ws.Range("H86").Value = 4900
$ws.Range("I86").Value = 4900
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4900
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3777
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4900
$ws.Range("I89").Value = 4900
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 24500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -18884
$ws.Range("N89").ClearContents()
$ws.Range("H102").Value = 39998.332
$ws.Range("I102").Value = 34998.5
$ws.Range("K102").Value = 34998.5
$ws.Range("M102").Value = -32564.5
$ws.Range("H108").Value = 47200
$ws.Range("J108").Value = 47200
$ws.Range("L108").Value = 47200
$ws.Range("N108").Value = -54880
$ws.Range("H109").Value = 24999
$ws.Range("J109").Value = 24999
$ws.Range("L109").Value = 24999
$ws.Range("N109").Value = -27079
$ws.Range("H112").Value = 80000
$ws.Range("J112").Value = 80000
$ws.Range("L112").Value = 80000
$ws.Range("N112").Value = -82954
$ws.Range("H132").Value = 10235.25
$ws.Range("J132").Value = 24457
$ws.Range("L132").Value = 73371
$ws.Range("N132").Value = -78431

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 22118.361
$ws.Range("I102").Value = 25010.932
$ws.Range("K102").Value = 25010.932
$ws.Range("M102").Value = -23388.932
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H128").Value = 34574.875
$ws.Range("J128").Value = 34574.875
$ws.Range("L128").Value = 34574.875
$ws.Range("N128").Value = -44534.875
$ws.Range("H132").Value = 3053.457
$ws.Range("J132").Value = 2367.2354
$ws.Range("L132").Value = 7101.706200000001
$ws.Range("N132").Value = -12161.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H82").Value = 1560.875
$ws.Range("I82").Value = 1337.0769
$ws.Range("J82").Value = 2530.6667
$ws.Range("K82").Value = 1337.0769
$ws.Range("L82").Value = 2530.6667
$ws.Range("M82").Value = -976.0769
$ws.Range("N82").Value = -3252.6667
$ws.Range("H85").Value = 1560.875
$ws.Range("I85").Value = 1337.0769
$ws.Range("J85").Value = 2530.6667
$ws.Range("K85").Value = 1337.0769
$ws.Range("L85").Value = 2530.6667
$ws.Range("M85").Value = -89.07690000000002
$ws.Range("N85").Value = -5026.6667
$ws.Range("H136").Value = 24398.146
$ws.Range("I136").Value = 2397.9167
$ws.Range("K136").Value = 7193.750100000001
$ws.Range("M136").Value = -4643.750100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 730.0714
$ws.Range("I107").Value = 838.55554
$ws.Range("K107").Value = 2515.66662
$ws.Range("M107").Value = -595.66662
$ws.Range("H132").Value = 2630.379
$ws.Range("I132").Value = 2148.228
$ws.Range("K132").Value = 6444.684
$ws.Range("M132").Value = -3914.684
$ws.Range("H136").Value = 3720.6182
$ws.Range("J136").Value = 4141.727
$ws.Range("L136").Value = 12425.181
$ws.Range("N136").Value = -17525.181
